$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.507663369178772
$ws.Range("B1").Value = 2.352208375930786
$ws.Range("C1").Value = 5.276742458343506
$ws.Range("D1").Value = 3.581472635269165
$ws.Range("E1").Value = 1.030079960823059
